$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @("利欧股份", "白银有色", "航天电子"),
    @("航天电子", "湖南白银", "白银有色"),
    @("白银有色", "利欧股份", "紫金矿业"),
    @("湖南白银", "航天电子", "浙文互联"),
    @("天奇股份", "湖南黄金", "利欧股份"),
    @("紫金矿业", "紫金矿业", "湖南白银"),
    @("蓝色光标", "金风科技", "通宇通讯"),
    @("浙文互联", "铜陵有色", "九鼎新材"),
    @("东方日升", "蓝色光标", "蓝色光标"),
    @("网宿科技", "东方日升", "铜陵有色"),
    @("铜陵有色", "隆基绿能", "航天发展"),
    @("金风科技", "浙文互联", "乾照光电"),
    @("湖南黄金", "天奇股份", "金风科技"),
    @("通宇通讯", "中国卫星", "洲际油气"),
    @("乾照光电", "乾照光电", "网宿科技"),
    @("中国卫星", "特变电工", "泓博医药"),
    @("九鼎新材", "通宇通讯", "红宝丽"),
    @("航天发展", "洛阳钼业", "天奇股份"),
    @("大位科技", "拓日新能", "洛阳钼业"),
    @("数据港", "航天发展", "中国卫星")
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $rowValues = $values[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
